$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row (253) down into the new rows,
# split into two ranges so columns H,I,J (score/result) and AB,AC (post-match PL odds)
# -- which do not apply to these not-yet-played fixtures -- are never populated.
$srcLeft = $ws.Range("A253:G253")
$srcRight = $ws.Range("K253:AA253")
$dstLeft = $ws.Range("A254:G262")
$dstRight = $ws.Range("K254:AA262")
$srcLeft.Copy($dstLeft)
$srcRight.Copy($dstRight)

# Row 254
$ws.Range("A254").Value = 252
$ws.Range("B254").Value = 6899154
$ws.Range("C254").Value = "Portugal Segunda Liga"
$ws.Range("D254").Value = "Portugal Segunda Liga"
$ws.Range("E254").Value = 45395.29166666666
$ws.Range("F254").Value = "CF Os Belenenses"
$ws.Range("G254").Value = "Academico Viseu"
$ws.Range("K254").Value = 3.25
$ws.Range("L254").Value = 3
$ws.Range("M254").Value = 2.1
$ws.Range("N254").Value = 3.3
$ws.Range("O254").Value = 3
$ws.Range("P254").Value = 2.05
$ws.Range("Q254").Value = 0.25
$ws.Range("R254").Value = 2.025
$ws.Range("S254").Value = 1.825
$ws.Range("T254").Value = 2.25
$ws.Range("U254").Value = 1.925
$ws.Range("V254").Value = 1.925
$ws.Range("W254").Value = 0
$ws.Range("X254").Value = 0
$ws.Range("Y254").Value = 0
$ws.Range("Z254").Value = 0
$ws.Range("AA254").Value = 0

# Row 255
$ws.Range("A255").Value = 253
$ws.Range("B255").Value = 6899155
$ws.Range("C255").Value = "Portugal Segunda Liga"
$ws.Range("D255").Value = "Portugal Segunda Liga"
$ws.Range("E255").Value = 45395.41666666666
$ws.Range("F255").Value = "UD Leiria"
$ws.Range("G255").Value = "Vilaverdense"
$ws.Range("K255").Value = 1.6
$ws.Range("L255").Value = 3.75
$ws.Range("M255").Value = 5
$ws.Range("N255").Value = 1.6
$ws.Range("O255").Value = 3.75
$ws.Range("P255").Value = 5
$ws.Range("Q255").Value = -0.75
$ws.Range("R255").Value = 1.8
$ws.Range("S255").Value = 2.05
$ws.Range("T255").Value = 2.5
$ws.Range("U255").Value = 1.9
$ws.Range("V255").Value = 1.95
$ws.Range("W255").Value = 0
$ws.Range("X255").Value = 0
$ws.Range("Y255").Value = 0
$ws.Range("Z255").Value = 0
$ws.Range("AA255").Value = 0

# Row 256
$ws.Range("A256").Value = 254
$ws.Range("B256").Value = 6893187
$ws.Range("C256").Value = "Portugal Segunda Liga"
$ws.Range("D256").Value = "Portugal Segunda Liga"
$ws.Range("E256").Value = 45395.47916666666
$ws.Range("F256").Value = "Tondela"
$ws.Range("G256").Value = "Penafiel"
$ws.Range("K256").Value = 1.9
$ws.Range("L256").Value = 3.25
$ws.Range("M256").Value = 4
$ws.Range("N256").Value = 1.85
$ws.Range("O256").Value = 3.3
$ws.Range("P256").Value = 4.2
$ws.Range("Q256").Value = -0.5
$ws.Range("R256").Value = 1.9
$ws.Range("S256").Value = 1.95
$ws.Range("T256").Value = 2.25
$ws.Range("U256").Value = 1.975
$ws.Range("V256").Value = 1.875
$ws.Range("W256").Value = 0
$ws.Range("X256").Value = 0
$ws.Range("Y256").Value = 0
$ws.Range("Z256").Value = 0
$ws.Range("AA256").Value = 0

# Row 257
$ws.Range("A257").Value = 255
$ws.Range("B257").Value = 6899262
$ws.Range("C257").Value = "Portugal Segunda Liga"
$ws.Range("D257").Value = "Portugal Segunda Liga"
$ws.Range("E257").Value = 45395.47916666666
$ws.Range("F257").Value = "Benfica B"
$ws.Range("G257").Value = "AVS"
$ws.Range("K257").Value = 2.75
$ws.Range("L257").Value = 3.1
$ws.Range("M257").Value = 2.5
$ws.Range("N257").Value = 2.8
$ws.Range("O257").Value = 3.1
$ws.Range("P257").Value = 2.45
$ws.Range("Q257").Value = 0
$ws.Range("R257").Value = 2.1
$ws.Range("S257").Value = 1.775
$ws.Range("T257").Value = 2.5
$ws.Range("U257").Value = 2.025
$ws.Range("V257").Value = 1.825
$ws.Range("W257").Value = 0
$ws.Range("X257").Value = 0
$ws.Range("Y257").Value = 0
$ws.Range("Z257").Value = 0
$ws.Range("AA257").Value = 0

# Row 258
$ws.Range("A258").Value = 256
$ws.Range("B258").Value = 6893188
$ws.Range("C258").Value = "Portugal Segunda Liga"
$ws.Range("D258").Value = "Portugal Segunda Liga"
$ws.Range("E258").Value = 45396.29166666666
$ws.Range("F258").Value = "Pacos Ferreira"
$ws.Range("G258").Value = "Nacional"
$ws.Range("K258").Value = 2.4
$ws.Range("L258").Value = 3.2
$ws.Range("M258").Value = 2.8
$ws.Range("N258").Value = 2.375
$ws.Range("O258").Value = 3.2
$ws.Range("P258").Value = 2.875
$ws.Range("Q258").Value = -0.25
$ws.Range("R258").Value = 2.1
$ws.Range("S258").Value = 1.775
$ws.Range("T258").Value = 2.5
$ws.Range("U258").Value = 1.975
$ws.Range("V258").Value = 1.875
$ws.Range("W258").Value = 0
$ws.Range("X258").Value = 0
$ws.Range("Y258").Value = 0
$ws.Range("Z258").Value = 0
$ws.Range("AA258").Value = 0

# Row 259
$ws.Range("A259").Value = 257
$ws.Range("B259").Value = 6893617
$ws.Range("C259").Value = "Portugal Segunda Liga"
$ws.Range("D259").Value = "Portugal Segunda Liga"
$ws.Range("E259").Value = 45396.41666666666
$ws.Range("F259").Value = "CD Mafra"
$ws.Range("G259").Value = "Feirense"
$ws.Range("K259").Value = 1.9
$ws.Range("L259").Value = 3.4
$ws.Range("M259").Value = 3.75
$ws.Range("N259").Value = 1.75
$ws.Range("O259").Value = 3.5
$ws.Range("P259").Value = 4.2
$ws.Range("Q259").Value = -0.75
$ws.Range("R259").Value = 2.05
$ws.Range("S259").Value = 1.8
$ws.Range("T259").Value = 2.5
$ws.Range("U259").Value = 1.975
$ws.Range("V259").Value = 1.875
$ws.Range("W259").Value = 0
$ws.Range("X259").Value = 0
$ws.Range("Y259").Value = 0
$ws.Range("Z259").Value = 0
$ws.Range("AA259").Value = 0

# Row 260
$ws.Range("A260").Value = 258
$ws.Range("B260").Value = 6893189
$ws.Range("C260").Value = "Portugal Segunda Liga"
$ws.Range("D260").Value = "Portugal Segunda Liga"
$ws.Range("E260").Value = 45396.47916666666
$ws.Range("F260").Value = "Leixoes"
$ws.Range("G260").Value = "SCU Torreense"
$ws.Range("K260").Value = 2.2
$ws.Range("L260").Value = 3
$ws.Range("M260").Value = 3.4
$ws.Range("N260").Value = 2.2
$ws.Range("O260").Value = 3
$ws.Range("P260").Value = 3.4
$ws.Range("Q260").Value = -0.25
$ws.Range("R260").Value = 1.9
$ws.Range("S260").Value = 1.95
$ws.Range("T260").Value = 2
$ws.Range("U260").Value = 1.8
$ws.Range("V260").Value = 2.05
$ws.Range("W260").Value = 0
$ws.Range("X260").Value = 0
$ws.Range("Y260").Value = 0
$ws.Range("Z260").Value = 0
$ws.Range("AA260").Value = 0

# Row 261
$ws.Range("A261").Value = 259
$ws.Range("B261").Value = 6893618
$ws.Range("C261").Value = "Portugal Segunda Liga"
$ws.Range("D261").Value = "Portugal Segunda Liga"
$ws.Range("E261").Value = 45396.47916666666
$ws.Range("F261").Value = "FC Porto B"
$ws.Range("G261").Value = "UD Oliveirense"
$ws.Range("K261").Value = 1.666
$ws.Range("L261").Value = 3.75
$ws.Range("M261").Value = 4.5
$ws.Range("N261").Value = 1.65
$ws.Range("O261").Value = 3.75
$ws.Range("P261").Value = 4.5
$ws.Range("Q261").Value = -0.75
$ws.Range("R261").Value = 1.875
$ws.Range("S261").Value = 1.975
$ws.Range("T261").Value = 2.75
$ws.Range("U261").Value = 1.875
$ws.Range("V261").Value = 1.975
$ws.Range("W261").Value = 0
$ws.Range("X261").Value = 0
$ws.Range("Y261").Value = 0
$ws.Range("Z261").Value = 0
$ws.Range("AA261").Value = 0

# Row 262
$ws.Range("A262").Value = 260
$ws.Range("B262").Value = 6899156
$ws.Range("C262").Value = "Portugal Segunda Liga"
$ws.Range("D262").Value = "Portugal Segunda Liga"
$ws.Range("E262").Value = 45396.6875
$ws.Range("F262").Value = "Maritimo"
$ws.Range("G262").Value = "Santa Clara"
$ws.Range("K262").Value = 2.6
$ws.Range("L262").Value = 3
$ws.Range("M262").Value = 2.75
$ws.Range("N262").Value = 2.55
$ws.Range("O262").Value = 3
$ws.Range("P262").Value = 2.8
$ws.Range("Q262").Value = 0
$ws.Range("R262").Value = 1.825
$ws.Range("S262").Value = 2.025
$ws.Range("T262").Value = 2
$ws.Range("U262").Value = 1.9
$ws.Range("V262").Value = 1.95
$ws.Range("W262").Value = 0
$ws.Range("X262").Value = 0
$ws.Range("Y262").Value = 0
$ws.Range("Z262").Value = 0
$ws.Range("AA262").Value = 0

Write-Host "Added rows 254-262"
